$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (Employee_ID 5456 / Dulce / Human Resource) - all rows below shift up
$ws.Rows.Item(2).Delete()

# Update selection to match the post-edit state (A2:C2 selected, active cell A2)
$ws.Range("A2:C2").Select()
